$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels
$ws.Range("E2").Value = "moyenne hantel"
$ws.Range("G2").Value = "E-T hantel"

# Update summary statistics (F1, H1, F2, H2)
$ws.Range("F1").Value = 1.685193628587822
$ws.Range("H1").Value = 0.7395101602815672
$ws.Range("F2").Value = 0.436885088002929
$ws.Range("H2").Value = 0.5338222578240785

# Update column A (rows 3-15)
$ws.Range("A3").Value = 1.847422334842916
$ws.Range("A4").Value = 2.280350850198276
$ws.Range("A5").Value = 1.958717916262161
$ws.Range("A6").Value = 1.675774202594786
$ws.Range("A7").Value = 1.893308053280127
$ws.Range("A8").Value = 2.12184406171194
$ws.Range("A9").Value = 1.788854381999832
$ws.Range("A10").Value = 1.629643428765334
$ws.Range("A11").Value = 2.252370025485892
$ws.Range("A13").Value = 2.100114365777487
$ws.Range("A14").Value = 1.969966809560311
$ws.Range("A15").Value = 2.07434436975044

# Update column B (only rows 3-7 and row 10 change, rest stay 0)
$ws.Range("B3").Value = 1.16521230467429
$ws.Range("B4").Value = 1.060954090791133
$ws.Range("B5").Value = 1.064737449050925
$ws.Range("B6").Value = 0.9540319654762321
$ws.Range("B7").Value = 1.157495537083123
$ws.Range("B10").Value = 0.7139598849653022
